$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-18 04:41:31"

# zh-cn sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-18 04:41:15"
$wsZhCn.Range("K2").Value = "2016-10-18 04:42:08"

# de-de sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-18 04:41:31"
$wsDeDe.Range("K2").Value = "2016-10-18 04:42:31"
